# Add a new row (row 91) of price data to the worksheet, following the
# same layout as the existing rows: Date (text), Original (CNY/mt),
# VAT Included (USD/mt), VAT Excluded (USD/mt), USD/CNY.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(91, 1).Value = "2024-11-01 00:00:00"
$ws.Cells.Item(91, 2).Value = 73850
$ws.Cells.Item(91, 3).Value = 10340.96
$ws.Cells.Item(91, 4).Value = 9151.299999999999
$ws.Cells.Item(91, 5).Value = 7.1237
